$wb = $excel.ActiveWorkbook

# --- GET_Tests sheet (was not the active tab; becomes the active tab) ---
$wsGet = $wb.Worksheets.Item("GET_Tests")

# Remove the stray L3 cell (value 750) from row 3
$wsGet.Range("L3").ClearContents()

# --- POST Tests sheet (was the active tab; loses that status) ---
$wsPost = $wb.Worksheets.Item("POST Tests")

# Remove the stray L1 cell (value 750) from row 1
$wsPost.Range("L1").ClearContents()

# Update POST Tests' own selection/scroll position before it stops being active
$wsPost.Range("L1").Select()

# Make GET_Tests the active sheet/tab (activeTab=0, tabSelected flips between sheets)
$wsGet.Activate()

# Restore GET_Tests' selection to L3 (matches the unchanged <selection> in the diff)
$wsGet.Range("L3").Select()
